# Refresh the FASS department course/requisites table with the
# Spring 2026 (1261) term data: a few course descriptions were revised
# (FASS222/223/224), two new courses were added (FASS100, FASS303), and
# the table was re-sorted by Course_Code.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Clear the previous data range (keep header) and the used range ----
$usedRange = $ws.UsedRange
$usedRange.ClearContents()

# ---- Header row ----
$headers = @("Course_Code", "Type", "Dept", "Course_Name", "Prerequisites", "CoRequisites", "Acad_Group", "Acad_Org", "Description")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# ---- Data rows, sorted alphabetically by Course_Code ----
$rows = @(
    @("FASS100", "Course", "FASS", "FASS 100 - STT-Introduction to Higher Education", "None", "None", "ARTS", "DEAN ARTS", ""),
    @("FASS201", "Course", "FASS", "FASS 201 - FASS Co-op Practicum I", "None", "None", "ARTS", "DEAN ARTS", ""),
    @("FASS222", "Course", "FASS", "FASS 222 - Academics in FASS", "None", "None", "ARTS", "DEAN ARTS", ""),
    @("FASS223", "Course", "FASS", "FASS 223 - Finding a Career in the Arts and Social Sciences", "None", "None", "ARTS", "DEAN ARTS", ""),
    @("FASS224", "Course", "FASS", "FASS 224 - Arts and Social Sciences in Everyday Life", "None", "None", "ARTS", "DEAN ARTS", ""),
    @("FASS225", "Course", "FASS", "FASS 225 - Special Topics in Basic Data Analysis in the Arts and Social Sciences", "None", "None", "ARTS", "DEAN ARTS", "REQ-15 units. Students with credit for FASS 208, FASS 210 or FASS 211 may not take this course for further credit unless a different topic is offered."),
    @("FASS301", "Course", "FASS", "FASS 301 - FASS Co-op Practicum II", "None", "None", "ARTS", "DEAN ARTS", ""),
    @("FASS302", "Course", "FASS", "FASS 302 - FASS Co-op Practicum III", "None", "None", "ARTS", "DEAN ARTS", ""),
    @("FASS303", "Course", "FASS", "FASS 303 - STT-Critical Information Literacy in the Humanities and Social Sciences", "None", "None", "ARTS", "DEAN ARTS", ""),
    @("FASS305", "Course", "FASS", "FASS 305 - STT-Centring Justice, Equity, Diversity, and Inclusion (JEDI) in Practical Life Settings", "None", "None", "ARTS", "DEAN ARTS", ""),
    @("FASS401", "Course", "FASS", "FASS 401 - FASS Co-op Practicum IV", "None", "None", "ARTS", "DEAN ARTS", ""),
    @("FASS402", "Course", "FASS", "FASS 402 - FASS Co-op Practicum V", "None", "None", "ARTS", "DEAN ARTS", ""),
    @("FASS499", "Course", "FASS", "FASS 499 - STT-FASS Dean's Leadership Seminar", "None", "None", "ARTS", "DEAN ARTS", ""),
    @("FASS805", "Course", "FASS", "FASS 805 - STT-Centring Justice, Equity, Diversity, and Inclusion (JEDI) in Practical Life Settings", "None", "None", "ARTS", "DEAN ARTS", "")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowData = $rows[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $value = $rowData[$c]
        if ($value -ne "") {
            $ws.Cells.Item($excelRow, $c + 1).Value = $value
        }
    }
}

# ---- Sheet view / selection tidy-up ----
$ws.Range("N5").Select()
